$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously-empty "想找的工作" (A) / "现在的状态" (B) entries
# for the newly added wishes (rows 55-75).

$ws.Range("A55").Value = "ai产品"
$ws.Range("B55").Value = "在找"

$ws.Range("A56").Value = "ai产品"
$ws.Range("B56").Value = "在找"

$ws.Range("A57").Value = "数据开发+ai"
$ws.Range("B57").Value = "失业了在找"

$ws.Range("A58").Value = "AI+JAVA"
$ws.Range("B58").Value = "在找"

$ws.Range("A59").Value = "Java实习"
$ws.Range("B59").Value = "在找"

$ws.Range("A60").Value = "前端工程师"
$ws.Range("B60").Value = "大专，前端工程师，空窗创业，gap一年半😁，找了内推看看机会`n保持学习，看下金九银十有没机会了"

$ws.Range("A61").Value = "AI产品"
$ws.Range("B61").Value = "在找"

$ws.Range("A62").Value = "Java秋招"
$ws.Range("B62").Value = "目前实习，在改简历背八股准备秋招！！！"

$ws.Range("A63").Value = "大厂的产品岗_26届秋招"
$ws.Range("B63").Value = "应届生正在进行中"

$ws.Range("A64").Value = "双休、Java"
$ws.Range("B64").Value = "25应届摆烂仔正在行动"

$ws.Range("A65").Value = "AI + JAVA "
$ws.Range("B65").Value = "在找"

$ws.Range("A66").Value = "嵌入式Ai"
$ws.Range("B66").Value = "在找"

$ws.Range("A67").Value = "全栈架构师 研发经理"
$ws.Range("B67").Value = "被裁再找,改好了简历"

$ws.Range("A68").Value = "AI产品"
$ws.Range("B68").Value = "在找"

$ws.Range("A69").Value = "AI应用开发"
$ws.Range("B69").Value = "在找"

$ws.Range("A70").Value = "网络安全"
$ws.Range("B70").Value = "在找"

$ws.Range("A71").Value = "java开发"
$ws.Range("B71").Value = "在找"

$ws.Range("A72").Value = "嵌入式/硬件测试"
$ws.Range("B72").Value = "25届毕业生在狂卷"

$ws.Range("A73").Value = "Java 中厂"
$ws.Range("B73").Value = "刚找到实习，现在边实习边秋招"

$ws.Range("A74").Value = "Java+AI"
$ws.Range("B74").Value = "在找"

$ws.Range("A75").Value = "AIGC相关后端或全栈"
$ws.Range("B75").Value = "开始找第一天"
